$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.413.01'
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = '  +3.60%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.795.32'
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = '  +4.15%  '

$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '337.29'
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = '  +1.61%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = '  -0.26%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3801'
$ws.Range("D7").Style = "Normal"

$ws.Range("E7").Value = '  +1.82%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3443'
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").Value = '  +1.85%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '48.15'
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Value = '  +0.00%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.204'
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = '  +2.08%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07515'
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").Value = '  +1.15%  '

$ws.Range("E12").Value = '  -0.26%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.07'
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = '  +10.39%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.488'
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = '  +1.57%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.793.82'
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = '  +4.21%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.078'
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = '  +0.76%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001099'
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = '  +2.55%  '

$ws.Range("E18").Value = '  -0.12%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '84.90'
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = '  +3.64%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.001'
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = '  -0.19%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.40'
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = '  +5.40%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.495'
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = '  +5.30%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.386.08'
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Value = '  +3.54%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.54'
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").Value = '  -1.22%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.448'
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").Value = '  -0.40%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.591'
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").Value = '  +8.80%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.495'
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").Value = '  +6.29%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '21.42'
$ws.Range("D28").Style = "Normal"

$ws.Range("E28").Value = '  +10.43%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '151.54'
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").Value = '  +0.57%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.994.49'
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").Value = '  +4.24%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '133.61'
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").Value = '  +1.83%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.066'
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").Value = '  -0.69%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.144'
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").Value = '  +3.26%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08717'
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").Value = '  +1.19%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '13.29'
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").Value = '  +4.41%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.682'
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").Value = '  -0.24%  '

$ws.Range("B37").Value = 'InternetComputer(DFINITY)'

$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.451'
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").Value = '  +1.89%  '

$ws.Range("B38").Value = 'TheSandbox'

$ws.Range("C38").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6909'
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").Value = '  +11.49%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.876'
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").Value = '  +5.93%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06356'
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").Value = '  +2.49%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2203'
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = '  +2.38%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.02341'
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = '  +0.44%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.277'
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = '  +4.68%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.47'
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = '  +1.90%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6465'
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = '  +7.63%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.000'
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Value = '  -0.21%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.862'
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Value = '  -0.76%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.120'
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").Value = '  +4.01%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '129.88'
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = '  +1.15%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07188'
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").Value = '  +0.31%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '79.23'
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").Value = '  +3.15%  '
